$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.768.91'
$ws.Range("E2").Value = '  +4.96%  '

$ws.Range("D3").Value = '3.967.42'
$ws.Range("E3").Value = '  +2.56%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +11.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.05%  '

$ws.Range("E7").Value = '  +1.77%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.747'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.29%  '

$ws.Range("E10").Value = '  +4.62%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.02'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.39%  '

$ws.Range("E12").Value = '  +4.04%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.81'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.56%  '

$ws.Range("D14").Value = '4.605.28'
$ws.Range("E14").Value = '  +2.46%  '

$ws.Range("D15").Value = '3.962.29'
$ws.Range("E15").Value = '  +1.95%  '

$ws.Range("E16").Value = '  +11.20%  '

$ws.Range("E17").Value = '  +4.62%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.37'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.35%  '

$ws.Range("E19").Value = '  +0.68%  '

$ws.Range("D20").Value = '72.479.29'
$ws.Range("E20").Value = '  +4.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '433.56'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.50%  '

$ws.Range("E22").Value = '  +15.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '95.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.50%  '

$ws.Range("E24").Value = '  +0.62%  '

$ws.Range("E25").Value = '  +3.62%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +23.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.54'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.90%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.93'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.06%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.36'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.86'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.97%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.00%  '

$ws.Range("E33").Value = '  +5.74%  '

$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '48.82'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.56%  '

$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '679.04'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.76%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '69.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.63%  '

$ws.Range("E37").Value = '  +3.94%  '

$ws.Range("D38").Value = '0.0₃0868'
$ws.Range("E38").Value = '  +11.91%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.14%  '

$ws.Range("E40").Value = '  +1.09%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.33'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.87%  '

$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("B44").Value = 'THORChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +15.66%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0486'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.56%  '

$ws.Range("E46").Value = '  +3.22%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.65'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.41'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.90%  '

$ws.Range("E50").Value = '  +8.08%  '

$ws.Range("D51").Value = '2.803.88'
$ws.Range("E51").Value = '  +12.87%  '
